$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 422 (Brocoli, Primera/Segunda pair),
# shifting all subsequent rows (old 422-447) down to 424-449.
$ws.Rows("422:423").Insert()

# Row 422: new "Primera" quality entry for date 44516 (2021-11-16)
$ws.Cells.Item(422, 1).Value = 3
$ws.Cells.Item(422, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(422, 3).Value = "Coquimbo"
$ws.Range("D422").Value = 44516
$ws.Range("D422").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(422, 5).Value = 5
$ws.Cells.Item(422, 6).Value = 100112023
$ws.Cells.Item(422, 7).Value = "Brócoli"
$ws.Cells.Item(422, 8).Value = "Sin especificar"
$ws.Cells.Item(422, 9).Value = "Primera"
$ws.Cells.Item(422, 10).Value = 2400
$ws.Cells.Item(422, 11).Value = 450
$ws.Cells.Item(422, 12).Value = 600
$ws.Cells.Item(422, 13).Value = 525
$ws.Cells.Item(422, 14).Value = "$/unidad"
$ws.Cells.Item(422, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(422, 16).Value = 525
$ws.Cells.Item(422, 17).Value = 1
$ws.Cells.Item(422, 18).Value = "Hortaliza"

# Row 423: new "Segunda" quality entry for the same date 44516
$ws.Cells.Item(423, 1).Value = 3
$ws.Cells.Item(423, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(423, 3).Value = "Coquimbo"
$ws.Range("D423").Value = 44516
$ws.Range("D423").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(423, 5).Value = 5
$ws.Cells.Item(423, 6).Value = 100112023
$ws.Cells.Item(423, 7).Value = "Brócoli"
$ws.Cells.Item(423, 8).Value = "Sin especificar"
$ws.Cells.Item(423, 9).Value = "Segunda"
$ws.Cells.Item(423, 10).Value = 1300
$ws.Cells.Item(423, 11).Value = 500
$ws.Cells.Item(423, 12).Value = 500
$ws.Cells.Item(423, 13).Value = 500
$ws.Cells.Item(423, 14).Value = "$/unidad"
$ws.Cells.Item(423, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(423, 16).Value = 500
$ws.Cells.Item(423, 17).Value = 1
$ws.Cells.Item(423, 18).Value = "Hortaliza"
